$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D to Text format temporarily to preserve string values that look numeric
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.902.44"
$ws.Range("E2").Value = "  +1.42%  "
$ws.Range("D3").Value = "1.674.97"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "332.36"
$ws.Range("E5").Value = "  +8.00%  "
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "0.3648"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").Value = "47.07"
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("D9").Value = "0.3238"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").Value = "1.142"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("D11").Value = "0.07125"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").Value = "6.089"
$ws.Range("E13").Value = "  +2.58%  "
$ws.Range("D14").Value = "19.66"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("D15").Value = "1.668.71"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("D16").Value = "6.646"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "0.06551"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "0.9993"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "78.78"
$ws.Range("E20").Value = "  +2.46%  "
$ws.Range("D21").Value = "15.85"
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "5.917"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "12.83"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").Value = "24.929.20"
$ws.Range("E24").Value = "  +1.57%  "
$ws.Range("D25").Value = "2.441"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "2.392"
$ws.Range("E26").Value = "  +2.79%  "
$ws.Range("D27").Value = "148.22"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").Value = "18.69"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "1.854.02"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "126.43"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "1.180"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("D32").Value = "4.097"
$ws.Range("E32").Value = "  +0.97%  "
$ws.Range("D33").Value = "5.790"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("D34").Value = "0.08472"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("D35").Value = "1.660"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "12.30"
$ws.Range("E36").Value = "  -0.82%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").Value = "0.06046"
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "0.02235"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("D40").Value = "1.227"
$ws.Range("E40").Value = "  +1.75%  "
$ws.Range("E41").Value = "  +1.27%  "
$ws.Range("D42").Value = "8.227"
$ws.Range("E42").Value = "  +0.40%  "
$ws.Range("D43").Value = "0.9995"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("D44").Value = "0.5957"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("E45").Value = "  +8.52%  "
$ws.Range("D46").Value = "3.851"
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").Value = "0.5719"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("D48").Value = "124.37"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").Value = "1.963"
$ws.Range("E49").Value = "  +0.77%  "
$ws.Range("D50").Value = "0.07006"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("E51").Value = "  +3.62%  "

# Restore default style (remove temporary text format) to match original styling
$ws.Range("D2:D51").Style = "Normal"

